$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "Bahndau Vs luong"
$ws.Range("A40").Value = "Tranformers"

$ws.Range("A41").Select()
